$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z103").Value = "Trailing space test. "
